$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.657.59"
$ws.Range("E2").Value = "  +1.38%  "

# Row 3
$ws.Range("D3").Value = "1.795.64"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.76"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("E6").Value = "  +2.04%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.94"
$ws.Range("E8").Value = "  +3.60%  "

# Row 9
$ws.Range("E9").Value = "  +2.22%  "

# Row 10
$ws.Range("E10").Value = "  +1.06%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("E11").Value = "  +0.33%  "

# Row 12
$ws.Range("D12").Value = "2.056.02"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.14"
$ws.Range("E13").Value = "  +0.98%  "

# Row 14
$ws.Range("D14").Value = "1.790.82"
$ws.Range("E14").Value = "  +0.30%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.638"
$ws.Range("E15").Value = "  +2.46%  "

# Row 16
$ws.Range("D16").Value = "34.587.80"
$ws.Range("E16").Value = "  +1.41%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.28"
$ws.Range("E17").Value = "  +2.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.94"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.41"
$ws.Range("E19").Value = "  +0.87%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0803"
$ws.Range("E20").Value = "  +3.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.32"
$ws.Range("E21").Value = "  +3.82%  "

# Row 22
$ws.Range("E22").Value = "  -0.32%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.17"
$ws.Range("E23").Value = "  +1.65%  "

# Row 24
$ws.Range("E24").Value = "  +1.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.24"
$ws.Range("E25").Value = "  +2.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.27"
$ws.Range("E26").Value = "  +1.13%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.58"
$ws.Range("E27").Value = "  +1.62%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.116"
$ws.Range("E28").Value = "  +2.64%  "

# Row 29
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.04"
$ws.Range("E30").Value = "  +11.46%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0523"
$ws.Range("E31").Value = "  +0.77%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.80"
$ws.Range("E32").Value = "  +3.08%  "

# Row 33
$ws.Range("E33").Value = "  +0.08%  "

# Row 34
$ws.Range("E34").Value = "  +2.25%  "

# Row 35
$ws.Range("D35").Value = "1.422.99"
$ws.Range("E35").Value = "  -1.62%  "

# Row 36
$ws.Range("E36").Value = "  +6.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.678"
$ws.Range("E37").Value = "  +3.78%  "

# Row 38
$ws.Range("E38").Value = "  +0.68%  "

# Row 39
$ws.Range("E39").Value = "  +1.88%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.56"
$ws.Range("E40").Value = "  +6.62%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.40"
$ws.Range("E41").Value = "  +0.87%  "

# Row 42
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.76"
$ws.Range("E42").Value = "  +2.70%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.933"
$ws.Range("E43").Value = "  +1.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.56"
$ws.Range("E44").Value = "  +0.41%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0525"
$ws.Range("E45").Value = "  +3.31%  "

# Row 46
$ws.Range("E46").Value = "  +0.10%  "

# Row 47
$ws.Range("E47").Value = "  +0.49%  "

# Row 48
$ws.Range("D48").Value = "1.956.74"
$ws.Range("E48").Value = "  +0.66%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.36"
$ws.Range("E49").Value = "  +0.44%  "

# Row 50
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.15%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  -4.44%  "
